$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (D) and volume-change (E) values.
# Column D holds plain-text price strings (e.g. "3.521.64" uses dots as
# thousand separators, not a real number) in the source data, so any value
# that Excel could otherwise auto-parse as a genuine number is forced back
# to text (NumberFormat "@") before being written, matching the original
# inline-string cell content.

$ws.Range("D2").Value = "66.876.15"
$ws.Range("E2").Value = "  -3.69%  "
$ws.Range("D3").Value = "3.524.13"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.17"
$ws.Range("E5").Value = "  -5.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.31"
$ws.Range("E6").Value = "  -4.77%  "
$ws.Range("D7").Value = "3.521.07"
$ws.Range("E7").Value = "  -4.09%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("E10").Value = "  -4.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.77"
$ws.Range("E11").Value = "  -4.16%  "
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").Value = "  -5.28%  "
$ws.Range("D14").Value = "4.114.79"
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").Value = "3.502.12"
$ws.Range("E16").Value = "  -3.73%  "
$ws.Range("D17").Value = "66.855.96"
$ws.Range("E17").Value = "  -3.64%  "
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.29"
$ws.Range("E19").Value = "  -3.77%  "
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.50"
$ws.Range("E21").Value = "  -5.06%  "
$ws.Range("E22").Value = "  -8.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.625"
$ws.Range("E23").Value = "  -3.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.32"
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "3.660.62"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -6.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.12"
$ws.Range("E29").Value = "  -9.84%  "
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.63"
$ws.Range("E33").Value = "  -4.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.156"
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("D35").Value = "3.515.41"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").Value = "  -5.13%  "
$ws.Range("E37").Value = "  -7.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.93"
$ws.Range("E38").Value = "  -6.18%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "173.91"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.11"
$ws.Range("E42").Value = "  -4.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.51"
$ws.Range("E43").Value = "  -5.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0854"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.889"
$ws.Range("E45").Value = "  -3.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.32"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.17"
$ws.Range("E47").Value = "  -5.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.52"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.51"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.02"
$ws.Range("E51").Value = "  -4.67%  "
